$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / 1h-volume table with the latest
# values pulled by the scheduled "Updated cryptos list" GitHub Actions
# job (also re-syncs two pairs of rows whose ranking order changed:
# WrappedEther/Chainlink and PancakeSwap/Filecoin).
#
# All written values are plain text (matching the source's inline-string
# cells). Values that would otherwise be auto-parsed by Excel as numbers
# (e.g. "0.999", "11.15") are written through a temporary Text ("@")
# number format so they stay text, then the cell style is restored to
# "Normal" so no formatting changes are left behind.

$ws.Range('D2').Value = '34.688.79'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '1.795.41'
$ws.Range('E3').Value = '  +0.56%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '226.83'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +0.49%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.558'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +2.06%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -0.16%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '32.93'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +3.57%  '
$ws.Range('E9').Value = '  +2.09%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.0696'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  +1.16%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.0950'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = '2.054.15'
$ws.Range('E12').Value = '  +0.56%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '11.15'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.804.42'
$ws.Range('E14').Value = '  +0.94%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '0.637'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +2.24%  '
$ws.Range('D16').Value = '34.595.41'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('E17').Value = '  +2.52%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '68.91'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +1.13%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '248.36'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').Value = '0.0₃0804'
$ws.Range('E20').Value = '  +3.35%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '11.30'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +3.31%  '
$ws.Range('E22').Value = '  -0.17%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '4.19'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +2.07%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '165.42'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +2.46%  '
$ws.Range('E25').Value = '  +0.19%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '7.28'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +1.21%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '16.58'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +1.65%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '0.117'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +2.64%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '1.01'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +0.55%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '4.18'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +15.29%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '3.83'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +3.60%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '1.24'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +0.00%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.0523'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +0.79%  '
$ws.Range('E34').Value = '  +2.05%  '
$ws.Range('D35').Value = '1.423.43'
$ws.Range('E35').Value = '  -1.72%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '2.60'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +6.85%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.675'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +3.24%  '
$ws.Range('E38').Value = '  +1.98%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.0193'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +0.72%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '85.58'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +6.56%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '2.39'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +0.75%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.935'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +1.30%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '2.76'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +2.56%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '13.60'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('E45').Value = '  +3.46%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '6.08'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('D48').Value = '1.955.02'
$ws.Range('E48').Value = '  +0.53%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '106.31'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('E51').Value = '  -4.86%  '
